$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 98 ---
$ws.Range("A98").Value = 45477.2916666667
$ws.Range("C98").Value = 6.21999979019165

# --- Add new row 99 ---
# Copy the date-format style from A98 onto A99 before writing the value,
# so A99 ends up using the same (existing) style index instead of a new one.
$ws.Range("A98").Copy()
$ws.Range("A99").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A99").Value = 45478.5537152778

$ws.Range("B99").Value = 1500
$ws.Range("C99").Value = 6.26000022888184
$ws.Range("D99").Value = 6.21999979019165
$ws.Range("E99").Value = 6.21999979019165
$ws.Range("F99").Value = 6.26000022888184

# G99 looks numeric ("6.26000022888184") but must stay stored as text
# (shared string), matching column G's existing "adj_close" text convention.
$ws.Range("G99").NumberFormat = "@"
$ws.Range("G99").Value = "6.26000022888184"
$ws.Range("G99").ClearFormats()

$ws.Range("H99").Value = "PAL.MI"
